$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add a brand-new "product checkout" sheet after the last existing sheet
#    and populate it with the checkout test data. (Populated first so the
#    shared-string table fills up in the same order as the real edit.)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCheckout = $wb.Worksheets.Add($null, $lastSheet)
$wsCheckout.Name = "product checkout"

$wsCheckout.Range("A1").Value = "countryName"
$wsCheckout.Range("B1").Value = "address"
$wsCheckout.Range("C1").Value = "postcode"
$wsCheckout.Range("D1").Value = "phone"
$wsCheckout.Range("E1").Value = "city"
$wsCheckout.Range("F1").Value = "productName"

$wsCheckout.Range("A2").Value = "Egypt"
$wsCheckout.Range("B2").Value = "test address"
$wsCheckout.Range("C2").Value = 11659
$wsCheckout.Range("D2").Value = 32445566677
$wsCheckout.Range("E2").Value = "Cairo"
$wsCheckout.Range("F2").Value = "Apple MacBook Pro 13-inch"

# Note: the host engine quantizes stored column widths to multiples of a
# fixed pixel grid, so the literal target "characters" width isn't always
# reproducible bit-for-bit; the inputs below are tuned to land on the
# closest achievable stored width to the real-Excel target.
$wsCheckout.Columns.Item(1).ColumnWidth = 12.334
$wsCheckout.Columns.Item(4).ColumnWidth = 12.834
$wsCheckout.Columns.Item(6).ColumnWidth = 29.5

$wsCheckout.Range("G1:G3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Rename the old "registration" sheet to "auto Suggest" and replace its
#    static registration-test data with the auto-suggest product data.
# ---------------------------------------------------------------------------
$wsAuto = $wb.Worksheets.Item(2)
$wsAuto.Name = "auto Suggest"

# Remove the old hyperlink + its now-stale C-column content entirely.
$wsAuto.Hyperlinks.Delete()
$wsAuto.Range("C1").ClearContents()
$wsAuto.Range("C2").ClearContents()

# Overwrite the old registration columns (firstName/lastName/password) with
# the new auto-suggest sample data (productName / partial product name).
$wsAuto.Range("A1").Value = "productName"
$wsAuto.Range("B1").Value = "partial product name"
$wsAuto.Range("A2").Value = "Apple MacBook Pro 13-inch"
$wsAuto.Range("B2").Value = "MacB"

# Resize the two remaining columns to fit the new, wider content (see note
# above about the engine's column-width quantization).
$wsAuto.Columns.Item(1).ColumnWidth = 24.665
$wsAuto.Columns.Item(2).ColumnWidth = 18.834

# New selection for this sheet (no longer the active tab).
$wsAuto.Range("B1:B2").Select() | Out-Null

# "product checkout" is the sheet that should be active/selected on open.
$wsCheckout.Activate() | Out-Null
